$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) ---------------------------------------------------
# A1: new company name (keeps its existing right-aligned / bordered style)
$ws.Range("A1").Value = "Grupo Serquímica"

# B1:D1: cleared back to an (empty) text value but a leading apostrophe is
# used so the cells stay text-typed (quote-prefixed) instead of turning
# into numeric/blank cells - this mirrors what Excel does when a text
# cell is emptied via the formula bar instead of Delete.
$ws.Range("B1:D1").Value = "'"

# E1:L1: the sheet now spans out to column L. These brand-new header
# cells are general-aligned (the default for new cells) and are also
# quote-prefixed empty text.
$ws.Range("E1:L1").HorizontalAlignment = 1
$ws.Range("E1:L1").Value = "'"

# --- Row 2 --------------------------------------------------------------
$ws.Range("A2").Value = "azship tecnologia"
$ws.Range("B2:D2").Value = ""
$ws.Range("E2:L2").HorizontalAlignment = 1

# --- Row 3 ----------------------------------------------------------------
$ws.Range("A3").Value = "'"
$ws.Range("E3:L3").HorizontalAlignment = 1

# --- Row 4 ----------------------------------------------------------------
$ws.Range("A4").Value = "'"
$ws.Range("B4:D4").Value = ""
$ws.Range("E4:L4").HorizontalAlignment = 1

# --- Column widths --------------------------------------------------------
$ws.Columns("B").ColumnWidth = 19.862142857142857
$ws.Columns("C").ColumnWidth = 40.57642857142857
$ws.Columns("E:L").ColumnWidth = 13.576428571428572
